$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136, shifting existing rows 136:208 down to 137:209
$ws.Rows(136).Insert()

# Populate the newly inserted row 136 with the new data record
$ws.Range("A136").Value = 10
$ws.Range("B136").Value = "Vega Modelo de Temuco"
$ws.Range("C136").Value = "La Araucanía"
$ws.Range("D136").Value = 44523
$ws.Range("E136").Value = 9
$ws.Range("F136").Value = 100112001
$ws.Range("G136").Value = "Berenjena"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 40
$ws.Range("K136").Value = 10000
$ws.Range("L136").Value = 10000
$ws.Range("M136").Value = 10000
$ws.Range("N136").Value = "$/caja 60 unidades"
$ws.Range("O136").Value = "Región de Arica y Parinacota"
$ws.Range("P136").Value = 167
$ws.Range("Q136").Value = 60
$ws.Range("R136").Value = "Hortaliza"
